# Replace the "ser: 122" blog entry with a new "ser: 125" entry.
# (commit message: "125,lesson from life of Nuh (as)")
#
# This cell (I8 on Sheet1) currently holds the shared string:
#   type: blog
#   width: 2
#   height: 1
#   ser: 122
#
# It should instead hold:
#   type: blog
#   width: 2
#   height: 1
#   ser: 125

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newline = [char]10
$newValue = "type: blog" + $newline + "width: 2" + $newline + "height: 1" + $newline + "ser: 125"

$ws.Range("I8").Value = $newValue
